$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$win.SplitRow = 4
$win.Split = $true
$win.Panes.Item(2).Activate()
$excel.Goto($ws.Range("A11"), $true)
